$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typo in header B1: "locacalizacion" -> "localizacion"
$ws.Range("B1").Value = "localizacion"

# Remove the stray location value in B2 (it was "18:13:14:12S")
$ws.Range("B2").ClearContents()

# Update the active selection to B2
$ws.Range("B2").Select()
